$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values per diff
$ws.Range("B3").Value = 2.9
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 27
$ws.Range("B5").Value = 0.85

# Update the active selection to match the diff (C5 -> C4)
$ws.Range("C4").Select()
